$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 62
$ws.Range("H62").Value2 = 1288.4
$ws.Range("I62").Value2 = 0
$ws.Range("J62").Value2 = 1288.4
$ws.Range("K62").Value2 = 0
$ws.Range("L62").Value2 = 1288.4
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value2 = -2536.4
# row 65
$ws.Range("H65").Value2 = 1288.4
$ws.Range("I65").Value2 = 0
$ws.Range("J65").Value2 = 1288.4
$ws.Range("K65").Value2 = 0
$ws.Range("L65").ClearContents()
$ws.Range("M65").Value2 = 6442
$ws.Range("N65").Value2 = -12682
# row 69
$ws.Range("H69").Value2 = 3431.2
$ws.Range("I69").Value2 = 3062.4
$ws.Range("J69").Value2 = 3800
$ws.Range("K69").Value2 = 9187.200000000001
$ws.Range("L69").Value2 = 11400
$ws.Range("M69").Value2 = -8313.200000000001
$ws.Range("N69").Value2 = -13148
# row 72
$ws.Range("H72").Value2 = 3431.2
$ws.Range("I72").Value2 = 3062.4
$ws.Range("J72").Value2 = 3800
$ws.Range("K72").Value2 = 27561.6
$ws.Range("L72").Value2 = 34200
$ws.Range("M72").Value2 = -23193.6
$ws.Range("N72").Value2 = -42936
# row 76
$ws.Range("H76").Value2 = 5327.1333
$ws.Range("I76").Value2 = 3612.875
$ws.Range("J76").Value2 = 7286.2856
$ws.Range("K76").Value2 = 3612.875
$ws.Range("L76").Value2 = 7286.2856
$ws.Range("M76").Value2 = -3297.875
$ws.Range("N76").Value2 = -7916.2856
# row 79
$ws.Range("H79").Value2 = 5327.1333
$ws.Range("I79").Value2 = 3612.875
$ws.Range("J79").Value2 = 7286.2856
$ws.Range("K79").Value2 = 3612.875
$ws.Range("L79").Value2 = 7286.2856
$ws.Range("M79").Value2 = -2520.875
$ws.Range("N79").Value2 = -9470.285599999999
# row 80
$ws.Range("H80").Value2 = 1684.8889
$ws.Range("I80").Value2 = 1984
$ws.Range("J80").Value2 = 1086.6666
$ws.Range("K80").Value2 = 5952
$ws.Range("L80").Value2 = 3259.9998
$ws.Range("M80").Value2 = -4954
$ws.Range("N80").Value2 = -5255.9998
# row 83
$ws.Range("H83").Value2 = 1684.8889
$ws.Range("I83").Value2 = 1984
$ws.Range("J83").Value2 = 1086.6666
$ws.Range("K83").Value2 = 17856
$ws.Range("L83").Value2 = 9779.999400000001
$ws.Range("M83").Value2 = -12864
$ws.Range("N83").Value2 = -19763.9994
# row 112
$ws.Range("H112").Value2 = 1071.3208
$ws.Range("J112").Value2 = 1107.6
$ws.Range("L112").Value2 = 3322.8
$ws.Range("N112").Value2 = -5538.799999999999
# row 138
$ws.Range("H138").Value2 = 6174389.5
$ws.Range("I138").Value2 = 8773184
$ws.Range("J138").Value2 = 2252.125
$ws.Range("K138").Value2 = 26319552
$ws.Range("L138").Value2 = 6756.375
$ws.Range("M138").Value2 = -26314412
$ws.Range("N138").Value2 = -17036.375
# row 141
$ws.Range("H141").Value2 = 1050.2325
$ws.Range("I141").Value2 = 1050.2325
$ws.Range("J141").Value2 = 0
$ws.Range("K141").Value2 = 3150.6975
$ws.Range("L141").Value2 = 0
$ws.Range("M141").Value2 = 2029.3025
$ws.Range("N141").ClearContents()

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value2 = 7302.4253
$ws.Range("I32").Value2 = 7345
$ws.Range("J32").Value2 = 7122.6665
$ws.Range("K32").Value2 = 7345
$ws.Range("L32").Value2 = 7122.6665
$ws.Range("M32").Value2 = -7058
$ws.Range("N32").Value2 = -7696.6665
# row 61
$ws.Range("H61").Value2 = 10205267
$ws.Range("I61").Value2 = 10870762
$ws.Range("J61").Value2 = 996.6667
$ws.Range("K61").Value2 = 10870762
$ws.Range("L61").Value2 = 996.6667
$ws.Range("M61").Value2 = -10870550
$ws.Range("N61").Value2 = -1420.6667
# row 136
$ws.Range("H136").Value2 = 10205267
$ws.Range("I136").Value2 = 10870762
$ws.Range("J136").Value2 = 996.6667
$ws.Range("K136").Value2 = 32612286
$ws.Range("L136").Value2 = 2990.0001
$ws.Range("M136").Value2 = -32609736
$ws.Range("N136").Value2 = -8090.0001

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 134
$ws.Range("H134").Value2 = 2577.6567
$ws.Range("I134").Value2 = 1859.2909
$ws.Range("K134").Value2 = 5577.8727
$ws.Range("M134").Value2 = -3042.8727

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value2 = 4832797.5
$ws.Range("I31").Value2 = 1544.5
$ws.Range("J31").Value2 = 22225310
$ws.Range("K31").Value2 = 1544.5
$ws.Range("L31").Value2 = 22225310
$ws.Range("M31").Value2 = -1249.5
$ws.Range("N31").Value2 = -22225900
# row 34
$ws.Range("H34").Value2 = 4832797.5
$ws.Range("I34").Value2 = 1544.5
$ws.Range("J34").Value2 = 22225310
$ws.Range("K34").Value2 = 1544.5
$ws.Range("L34").Value2 = 22225310
$ws.Range("M34").Value2 = -1342.5
$ws.Range("N34").Value2 = -22225714
# row 58
$ws.Range("H58").Value2 = 1220.525
$ws.Range("I58").Value2 = 570.3333
$ws.Range("J58").Value2 = 4285.7144
$ws.Range("K58").Value2 = 570.3333
$ws.Range("L58").Value2 = 4285.7144
$ws.Range("M58").Value2 = -367.3333
$ws.Range("N58").Value2 = -4691.7144
# row 132
$ws.Range("H132").Value2 = 5883580
$ws.Range("I132").Value2 = 7813575
$ws.Range("J132").Value2 = 1690.0476
$ws.Range("K132").Value2 = 23440725
$ws.Range("L132").Value2 = 5070.142800000001
$ws.Range("M132").Value2 = -23438195
$ws.Range("N132").Value2 = -10130.1428
# row 134
$ws.Range("H134").Value2 = 1129.1351
$ws.Range("I134").Value2 = 1134.7192
$ws.Range("J134").Value2 = 1110.4117
$ws.Range("K134").Value2 = 3404.1576
$ws.Range("L134").Value2 = 3331.2351
$ws.Range("M134").Value2 = -869.1576
$ws.Range("N134").Value2 = -8401.2351
# row 136
$ws.Range("H136").Value2 = 1220.525
$ws.Range("I136").Value2 = 570.3333
$ws.Range("J136").Value2 = 4285.7144
$ws.Range("K136").Value2 = 1710.9999
$ws.Range("L136").Value2 = 12857.1432
$ws.Range("M136").Value2 = 839.0001
$ws.Range("N136").Value2 = -17957.1432

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 5
$ws.Range("H5").Value2 = 884
$ws.Range("I5").Value2 = 339.9091
$ws.Range("K5").Value2 = 1019.7273
$ws.Range("M5").Value2 = -907.7273
# row 98
$ws.Range("H98").Value2 = 300
$ws.Range("I98").Value2 = 300
$ws.Range("J98").Value2 = 0
$ws.Range("K98").Value2 = 900
$ws.Range("L98").Value2 = 0
$ws.Range("M98").Value2 = 598
$ws.Range("N98").ClearContents()
# row 132
$ws.Range("H132").Value2 = 3830
$ws.Range("I132").Value2 = 745
$ws.Range("J132").Value2 = 10000
$ws.Range("K132").Value2 = 6705
$ws.Range("L132").Value2 = 90000
$ws.Range("M132").Value2 = -4175
$ws.Range("N132").Value2 = -95060
# row 135
$ws.Range("H135").Value2 = 884
$ws.Range("I135").Value2 = 339.9091
$ws.Range("K135").Value2 = 3059.1819
$ws.Range("M135").Value2 = -524.1819

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 132
$ws.Range("H132").Value2 = 2951.2222
$ws.Range("I132").Value2 = 2120.5813
$ws.Range("J132").Value2 = 4737.1
$ws.Range("K132").Value2 = 6361.743899999999
$ws.Range("L132").Value2 = 14211.3
$ws.Range("M132").Value2 = -3831.743899999999
$ws.Range("N132").Value2 = -19271.3

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 40
$ws.Range("H40").Value2 = 12705.714
$ws.Range("I40").Value2 = 13988
$ws.Range("J40").Value2 = 9500
$ws.Range("K40").Value2 = 13988
$ws.Range("L40").Value2 = 9500
$ws.Range("M40").Value2 = -13852
$ws.Range("N40").Value2 = -9772
# row 132
$ws.Range("H132").Value2 = 5687289.5
$ws.Range("I132").Value2 = 3503.5
$ws.Range("J132").Value2 = 19240934
$ws.Range("K132").Value2 = 10510.5
$ws.Range("L132").Value2 = 57722802
$ws.Range("M132").Value2 = -7980.5
$ws.Range("N132").Value2 = -57727862

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 136
$ws.Range("H136").Value2 = 1141
$ws.Range("I136").Value2 = 667.2286
$ws.Range("J136").Value2 = 6668.3335
$ws.Range("K136").Value2 = 2001.6858
$ws.Range("L136").Value2 = 20005.0005
$ws.Range("M136").Value2 = 548.3141999999998
$ws.Range("N136").Value2 = -25105.0005
